$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row above the current row 27 (Ecuador / 2023-11-10),
# shifting all subsequent rows down by one (old row 99 -> new row 100).
$ws.Rows("27:27").Insert()

# Populate the new row 27 with the new weekly record (Peru / 2023-11-20).
$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 45250
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108007
$ws.Range("J27").Value = "Coco"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 30000
$ws.Range("O27").Value = 30000
$ws.Range("P27").Value = 30000
$ws.Range("Q27").Value = "$/malla 20 unidades"
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 1500
$ws.Range("T27").Value = 20
